$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(99, 8).Value = 1262481.5  # H99
$ws.Cells.Item(99, 9).Value = 1640616.9  # I99
$ws.Cells.Item(99, 10).Value = 2030  # J99
$ws.Cells.Item(99, 11).Value = 4921850.699999999  # K99
$ws.Cells.Item(99, 12).Value = 6090  # L99
$ws.Cells.Item(99, 13).Value = -4920352.699999999  # M99
$ws.Cells.Item(99, 14).Value = -9086  # N99

$ws.Cells.Item(132, 8).Value = 2131010  # H132
$ws.Cells.Item(132, 9).Value = 3216  # I132
$ws.Cells.Item(132, 11).Value = 9648  # K132
$ws.Cells.Item(132, 13).Value = -7118  # M132

$ws.Cells.Item(135, 8).Value = 7185  # H135
$ws.Cells.Item(135, 9).Value = 8085.4375  # I135
$ws.Cells.Item(135, 11).Value = 72768.9375  # K135
$ws.Cells.Item(135, 13).Value = -70233.9375  # M135

$ws.Cells.Item(136, 14).ClearContents()  # N136
$ws.Cells.Item(136, 8).Value = 0  # H136
$ws.Cells.Item(136, 10).Value = 0  # J136
$ws.Cells.Item(136, 12).Value = 0  # L136

$ws.Cells.Item(137, 8).Value = 9903.23  # H137
$ws.Cells.Item(137, 9).Value = 19037.166  # I137
$ws.Cells.Item(137, 10).Value = 2074.1428  # J137
$ws.Cells.Item(137, 11).Value = 57111.49800000001  # K137
$ws.Cells.Item(137, 12).Value = 6222.428400000001  # L137
$ws.Cells.Item(137, 13).Value = -54561.49800000001  # M137
$ws.Cells.Item(137, 14).Value = -11322.4284  # N137

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6386.6895  # H32
$ws.Cells.Item(32, 9).Value = 6347  # I32
$ws.Cells.Item(32, 10).Value = 7498  # J32
$ws.Cells.Item(32, 11).Value = 6347  # K32
$ws.Cells.Item(32, 12).Value = 7498  # L32
$ws.Cells.Item(32, 13).Value = -6060  # M32
$ws.Cells.Item(32, 14).Value = -8072  # N32

$ws.Cells.Item(61, 8).Value = 6491.079  # H61
$ws.Cells.Item(61, 9).Value = 7024.1035  # I61
$ws.Cells.Item(61, 11).Value = 7024.1035  # K61
$ws.Cells.Item(61, 13).Value = -6812.1035  # M61

$ws.Cells.Item(74, 8).Value = 4742.697  # H74
$ws.Cells.Item(74, 9).Value = 4944.7407  # I74
$ws.Cells.Item(74, 11).Value = 4944.7407  # K74
$ws.Cells.Item(74, 13).Value = -4070.7407  # M74

$ws.Cells.Item(77, 8).Value = 4742.697  # H77
$ws.Cells.Item(77, 9).Value = 4944.7407  # I77
$ws.Cells.Item(77, 11).Value = 24723.7035  # K77
$ws.Cells.Item(77, 13).Value = -20355.7035  # M77

$ws.Cells.Item(97, 8).Value = 5717605.5  # H97
$ws.Cells.Item(97, 9).Value = 4974.0454  # I97
$ws.Cells.Item(97, 10).Value = 15385136  # J97
$ws.Cells.Item(97, 11).Value = 4974.0454  # K97
$ws.Cells.Item(97, 12).Value = 15385136  # L97
$ws.Cells.Item(97, 13).Value = -4478.0454  # M97
$ws.Cells.Item(97, 14).Value = -15386128  # N97

$ws.Cells.Item(122, 8).Value = 772832  # H122
$ws.Cells.Item(122, 9).Value = 3004.9644  # I122
$ws.Cells.Item(122, 11).Value = 9014.893199999999  # K122
$ws.Cells.Item(122, 13).Value = -6564.893199999999  # M122

$ws.Cells.Item(123, 8).Value = 129999  # H123
$ws.Cells.Item(123, 10).Value = 129999  # J123
$ws.Cells.Item(123, 12).Value = 129999  # L123
$ws.Cells.Item(123, 14).Value = -139799  # N123

$ws.Cells.Item(132, 8).Value = 2632.4792  # H132
$ws.Cells.Item(132, 9).Value = 2287.8484  # I132
$ws.Cells.Item(132, 11).Value = 6863.5452  # K132
$ws.Cells.Item(132, 13).Value = -4333.5452  # M132

$ws.Cells.Item(136, 8).Value = 6491.079  # H136
$ws.Cells.Item(136, 9).Value = 7024.1035  # I136
$ws.Cells.Item(136, 11).Value = 21072.3105  # K136
$ws.Cells.Item(136, 13).Value = -18522.3105  # M136

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 147698.86  # H105
$ws.Cells.Item(105, 9).Value = 168982  # I105
$ws.Cells.Item(105, 10).Value = 20000  # J105
$ws.Cells.Item(105, 11).Value = 168982  # K105
$ws.Cells.Item(105, 12).Value = 20000  # L105
$ws.Cells.Item(105, 13).Value = -167235  # M105
$ws.Cells.Item(105, 14).Value = -23494  # N105

$ws.Cells.Item(107, 8).Value = 1816.2  # H107
$ws.Cells.Item(107, 9).Value = 1982.9231  # I107
$ws.Cells.Item(107, 10).Value = 732.5  # J107
$ws.Cells.Item(107, 11).Value = 1982.9231  # K107
$ws.Cells.Item(107, 12).Value = 732.5  # L107
$ws.Cells.Item(107, 13).Value = -62.92309999999998  # M107
$ws.Cells.Item(107, 14).Value = -4572.5  # N107

$ws.Cells.Item(134, 8).Value = 9289.117  # H134
$ws.Cells.Item(134, 9).Value = 9927.666999999999  # I134
$ws.Cells.Item(134, 11).Value = 29783.001  # K134
$ws.Cells.Item(134, 13).Value = -27248.001  # M134

$ws.Cells.Item(140, 8).Value = 84146.336  # H140
$ws.Cells.Item(140, 10).Value = 84146.336  # J140
$ws.Cells.Item(140, 12).Value = 84146.336  # L140
$ws.Cells.Item(140, 14).Value = -94506.336  # N140

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 13).ClearContents()  # M31
$ws.Cells.Item(31, 8).Value = 6281.2856  # H31
$ws.Cells.Item(31, 9).Value = 0  # I31
$ws.Cells.Item(31, 10).Value = 6281.2856  # J31
$ws.Cells.Item(31, 11).Value = 0  # K31
$ws.Cells.Item(31, 12).Value = 6281.2856  # L31
$ws.Cells.Item(31, 14).Value = -6871.2856  # N31

$ws.Cells.Item(34, 13).ClearContents()  # M34
$ws.Cells.Item(34, 8).Value = 6281.2856  # H34
$ws.Cells.Item(34, 9).Value = 0  # I34
$ws.Cells.Item(34, 10).Value = 6281.2856  # J34
$ws.Cells.Item(34, 11).Value = 0  # K34
$ws.Cells.Item(34, 12).Value = 6281.2856  # L34
$ws.Cells.Item(34, 14).Value = -6685.2856  # N34

$ws.Cells.Item(58, 8).Value = 2338.8125  # H58
$ws.Cells.Item(58, 9).Value = 2309  # I58
$ws.Cells.Item(58, 11).Value = 2309  # K58
$ws.Cells.Item(58, 13).Value = -2106  # M58

$ws.Cells.Item(109, 13).ClearContents()  # M109
$ws.Cells.Item(109, 8).Value = 53000  # H109
$ws.Cells.Item(109, 9).Value = 0  # I109
$ws.Cells.Item(109, 10).Value = 53000  # J109
$ws.Cells.Item(109, 11).Value = 0  # K109
$ws.Cells.Item(109, 12).Value = 53000  # L109
$ws.Cells.Item(109, 14).Value = -55080  # N109

$ws.Cells.Item(122, 8).Value = 7891.1333  # H122
$ws.Cells.Item(122, 9).Value = 6580.6816  # I122
$ws.Cells.Item(122, 11).Value = 19742.0448  # K122
$ws.Cells.Item(122, 13).Value = -17292.0448  # M122

$ws.Cells.Item(132, 8).Value = 2131.4546  # H132
$ws.Cells.Item(132, 9).Value = 2074.625  # I132
$ws.Cells.Item(132, 10).Value = 2283  # J132
$ws.Cells.Item(132, 11).Value = 6223.875  # K132
$ws.Cells.Item(132, 12).Value = 6849  # L132
$ws.Cells.Item(132, 13).Value = -3693.875  # M132
$ws.Cells.Item(132, 14).Value = -11909  # N132

$ws.Cells.Item(134, 8).Value = 7491.5454  # H134
$ws.Cells.Item(134, 9).Value = 9332.25  # I134
$ws.Cells.Item(134, 11).Value = 27996.75  # K134
$ws.Cells.Item(134, 13).Value = -25461.75  # M134

$ws.Cells.Item(136, 8).Value = 2338.8125  # H136
$ws.Cells.Item(136, 9).Value = 2309  # I136
$ws.Cells.Item(136, 11).Value = 6927  # K136
$ws.Cells.Item(136, 13).Value = -4377  # M136

$ws.Cells.Item(141, 8).Value = 297959.28  # H141
$ws.Cells.Item(141, 10).Value = 390143.3  # J141
$ws.Cells.Item(141, 12).Value = 390143.3  # L141
$ws.Cells.Item(141, 14).Value = -400503.3  # N141

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(132, 8).Value = 24043.863  # H132
$ws.Cells.Item(132, 9).Value = 1185  # I132
$ws.Cells.Item(132, 10).Value = 30767.059  # J132
$ws.Cells.Item(132, 11).Value = 10665  # K132
$ws.Cells.Item(132, 12).Value = 276903.531  # L132
$ws.Cells.Item(132, 13).Value = -8135  # M132
$ws.Cells.Item(132, 14).Value = -281963.531  # N132

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(47, 14).ClearContents()  # N47
$ws.Cells.Item(47, 8).Value = 0  # H47
$ws.Cells.Item(47, 10).Value = 0  # J47
$ws.Cells.Item(47, 12).Value = 0  # L47

$ws.Cells.Item(70, 8).Value = 8000.0435  # H70
$ws.Cells.Item(70, 9).Value = 6646.294  # I70
$ws.Cells.Item(70, 11).Value = 6646.294  # K70
$ws.Cells.Item(70, 13).Value = -6376.294  # M70

$ws.Cells.Item(73, 8).Value = 8000.0435  # H73
$ws.Cells.Item(73, 9).Value = 6646.294  # I73
$ws.Cells.Item(73, 11).Value = 6646.294  # K73
$ws.Cells.Item(73, 13).Value = -5710.294  # M73

$ws.Cells.Item(102, 8).Value = 6898.0386  # H102
$ws.Cells.Item(102, 9).Value = 8672.647000000001  # I102
$ws.Cells.Item(102, 11).Value = 8672.647000000001  # K102
$ws.Cells.Item(102, 13).Value = -7050.647000000001  # M102

$ws.Cells.Item(122, 8).Value = 9949.434999999999  # H122
$ws.Cells.Item(122, 9).Value = 7145.385  # I122
$ws.Cells.Item(122, 10).Value = 13594.7  # J122
$ws.Cells.Item(122, 11).Value = 21436.155  # K122
$ws.Cells.Item(122, 12).Value = 40784.10000000001  # L122
$ws.Cells.Item(122, 13).Value = -18986.155  # M122
$ws.Cells.Item(122, 14).Value = -45684.10000000001  # N122

$ws.Cells.Item(132, 8).Value = 3962.2444  # H132
$ws.Cells.Item(132, 9).Value = 4142.5557  # I132
$ws.Cells.Item(132, 11).Value = 12427.6671  # K132
$ws.Cells.Item(132, 13).Value = -9897.667099999999  # M132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 4609.067  # H61
$ws.Cells.Item(61, 9).Value = 830.5454999999999  # I61
$ws.Cells.Item(61, 10).Value = 15000  # J61
$ws.Cells.Item(61, 11).Value = 830.5454999999999  # K61
$ws.Cells.Item(61, 12).Value = 15000  # L61
$ws.Cells.Item(61, 13).Value = -628.5454999999999  # M61
$ws.Cells.Item(61, 14).Value = -15404  # N61

$ws.Cells.Item(64, 8).Value = 59000  # H64
$ws.Cells.Item(64, 10).Value = 59000  # J64
$ws.Cells.Item(64, 12).Value = 59000  # L64
$ws.Cells.Item(64, 14).Value = -59450  # N64

$ws.Cells.Item(67, 8).Value = 59000  # H67
$ws.Cells.Item(67, 10).Value = 59000  # J67
$ws.Cells.Item(67, 12).Value = 59000  # L67
$ws.Cells.Item(67, 14).Value = -60560  # N67

$ws.Cells.Item(113, 8).Value = 4609.067  # H113
$ws.Cells.Item(113, 9).Value = 830.5454999999999  # I113
$ws.Cells.Item(113, 10).Value = 15000  # J113
$ws.Cells.Item(113, 11).Value = 830.5454999999999  # K113
$ws.Cells.Item(113, 12).Value = 15000  # L113
$ws.Cells.Item(113, 13).Value = 1339.4545  # M113
$ws.Cells.Item(113, 14).Value = -19340  # N113

$ws.Cells.Item(136, 8).Value = 4642.9565  # H136
$ws.Cells.Item(136, 9).Value = 3351.7  # I136
$ws.Cells.Item(136, 11).Value = 10055.1  # K136
$ws.Cells.Item(136, 13).Value = -7505.099999999999  # M136

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(42, 8).Value = 13000  # H42
$ws.Cells.Item(42, 10).Value = 13000  # J42
$ws.Cells.Item(42, 12).Value = 13000  # L42
$ws.Cells.Item(42, 14).Value = -13756  # N42

$ws.Cells.Item(136, 8).Value = 561231.5600000001  # H136
$ws.Cells.Item(136, 9).Value = 677760.2  # I136
$ws.Cells.Item(136, 11).Value = 2033280.6  # K136
$ws.Cells.Item(136, 13).Value = -2030730.6  # M136
